# Auto-generated edit script applying cryptos.xlsx price/volume update (2024-11-21 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '96.729.64'
Set-TextCell $ws.Range("E2") '  +4.81%  '

Set-TextCell $ws.Range("D3") '3.128.21'
Set-TextCell $ws.Range("E3") '  +0.94%  '

Set-TextCell $ws.Range("E4") '  +0.00%  '

Set-TextCell $ws.Range("D5") '240.60'
Set-TextCell $ws.Range("E5") '  +3.27%  '

Set-TextCell $ws.Range("D6") '609.27'
Set-TextCell $ws.Range("E6") '  -0.64%  '

Set-TextCell $ws.Range("D7") '1.11'
Set-TextCell $ws.Range("E7") '  +2.15%  '

Set-TextCell $ws.Range("D8") '0.388'
Set-TextCell $ws.Range("E8") '  +0.40%  '

Set-TextCell $ws.Range("E9") '  +0.02%  '

Set-TextCell $ws.Range("B10") 'Cardano'
Set-TextCell $ws.Range("C10") 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws.Range("D10") '0.799'
Set-TextCell $ws.Range("E10") '  +2.24%  '

Set-TextCell $ws.Range("B11") 'LidoStakedEther'
Set-TextCell $ws.Range("C11") 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextCell $ws.Range("D11") '3.119.48'
Set-TextCell $ws.Range("E11") '  +0.79%  '

Set-TextCell $ws.Range("E12") '  +0.07%  '

Set-TextCell $ws.Range("D13") '95.878.85'
Set-TextCell $ws.Range("E13") '  +4.17%  '

Set-TextCell $ws.Range("D14") '0.0000242'
Set-TextCell $ws.Range("E14") '  -0.50%  '

Set-TextCell $ws.Range("D15") '34.30'
Set-TextCell $ws.Range("E15") '  +1.52%  '

Set-TextCell $ws.Range("D16") '5.36'
Set-TextCell $ws.Range("E16") '  -0.71%  '

Set-TextCell $ws.Range("D17") '3.698.58'
Set-TextCell $ws.Range("E17") '  +0.60%  '

Set-TextCell $ws.Range("D18") '3.105.70'
Set-TextCell $ws.Range("E18") '  +0.74%  '

Set-TextCell $ws.Range("D19") '3.60'
Set-TextCell $ws.Range("E19") '  -5.26%  '

Set-TextCell $ws.Range("B20") 'BitcoinCash'
Set-TextCell $ws.Range("C20") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Range("D20") '482.24'
Set-TextCell $ws.Range("E20") '  +10.53%  '

Set-TextCell $ws.Range("B21") 'Chainlink'
Set-TextCell $ws.Range("C21") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws.Range("D21") '14.52'
Set-TextCell $ws.Range("E21") '  +0.94%  '

Set-TextCell $ws.Range("D22") '5.74'
Set-TextCell $ws.Range("E22") '  +0.18%  '

Set-TextCell $ws.Range("D23") '0.0000196'
Set-TextCell $ws.Range("E23") '  -1.10%  '

Set-TextCell $ws.Range("D24") '8.83'
Set-TextCell $ws.Range("E24") '  -2.93%  '

Set-TextCell $ws.Range("D25") '5.56'
Set-TextCell $ws.Range("E25") '  +0.15%  '

Set-TextCell $ws.Range("D26") '85.49'
Set-TextCell $ws.Range("E26") '  +0.26%  '

Set-TextCell $ws.Range("D27") '11.78'
Set-TextCell $ws.Range("E27") '  +3.77%  '

Set-TextCell $ws.Range("D28") '3.264.40'
Set-TextCell $ws.Range("E28") '  +0.16%  '

Set-TextCell $ws.Range("E29") '  +0.05%  '

Set-TextCell $ws.Range("D30") '0.239'
Set-TextCell $ws.Range("E30") '  +1.69%  '

Set-TextCell $ws.Range("D31") '0.177'
Set-TextCell $ws.Range("E31") '  -0.74%  '

Set-TextCell $ws.Range("D32") '0.126'
Set-TextCell $ws.Range("E32") '  +2.73%  '

Set-TextCell $ws.Range("D34") '9.12'
Set-TextCell $ws.Range("E34") '  -0.07%  '

Set-TextCell $ws.Range("D35") '26.33'
Set-TextCell $ws.Range("E35") '  +3.00%  '

Set-TextCell $ws.Range("D36") '7.47'
Set-TextCell $ws.Range("E36") '  -6.96%  '

Set-TextCell $ws.Range("D37") '0.151'
Set-TextCell $ws.Range("E37") '  -2.33%  '

Set-TextCell $ws.Range("D38") '495.34'
Set-TextCell $ws.Range("E38") '  +6.41%  '

Set-TextCell $ws.Range("E39") '  -0.62%  '

Set-TextCell $ws.Range("D40") '24.18'
Set-TextCell $ws.Range("E40") '  +1.33%  '

Set-TextCell $ws.Range("D41") '0.441'
Set-TextCell $ws.Range("E41") '  +1.15%  '

Set-TextCell $ws.Range("B42") 'Fetch.AI'
Set-TextCell $ws.Range("C42") 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws.Range("D42") '1.24'
Set-TextCell $ws.Range("E42") '  -1.94%  '

Set-TextCell $ws.Range("B43") 'MantraDAO'
Set-TextCell $ws.Range("C43") 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextCell $ws.Range("D43") '3.64'
Set-TextCell $ws.Range("E43") '  -6.28%  '

Set-TextCell $ws.Range("B44") 'dogwifhat'
Set-TextCell $ws.Range("C44") 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell $ws.Range("D44") '3.23'
Set-TextCell $ws.Range("E44") '  -1.01%  '

Set-TextCell $ws.Range("B45") 'USDe'
Set-TextCell $ws.Range("C45") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell $ws.Range("D45") '1.00'
Set-TextCell $ws.Range("E45") '  -0.14%  '

Set-TextCell $ws.Range("D46") '161.87'
Set-TextCell $ws.Range("E46") '  +1.79%  '

Set-TextCell $ws.Range("D47") '0.702'
Set-TextCell $ws.Range("E47") '  +3.21%  '

Set-TextCell $ws.Range("D48") '1.91'
Set-TextCell $ws.Range("E48") '  +4.21%  '

Set-TextCell $ws.Range("D49") '44.04'
Set-TextCell $ws.Range("E49") '  +0.69%  '

Set-TextCell $ws.Range("B50") 'VeChain'
Set-TextCell $ws.Range("C50") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D50") '0.0323'
Set-TextCell $ws.Range("E50") '  -0.04%  '

Set-TextCell $ws.Range("D51") '4.37'
Set-TextCell $ws.Range("E51") '  +1.36%  '
